$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.511.59"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "2.288.74"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'546.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'131.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.573"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").Value = "2.287.76"
$ws.Range("E9").Value = "  -5.43%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.335"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.06%  "
$ws.Range("D14").Value = "'23.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("D15").Value = "2.693.42"
$ws.Range("E15").Value = "  -5.56%  "
$ws.Range("D16").Value = "58.422.52"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").Value = "2.297.77"
$ws.Range("E18").Value = "  -5.27%  "
$ws.Range("D19").Value = "'10.64"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.59%  "
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").Value = "'315.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "'6.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D24").Value = "'62.81"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("E25").Value = "  -4.09%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'8.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = "  -5.02%  "
$ws.Range("D29").Value = "'1.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'170.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -6.24%  "
$ws.Range("E32").Value = "  -5.01%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'17.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("D39").Value = "'3.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.97%  "
$ws.Range("D40").Value = "'37.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").Value = "'1.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("D42").Value = "'298.83"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.98%  "
$ws.Range("D43").Value = "'140.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("D44").Value = "'3.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.72%  "
$ws.Range("D45").Value = "'0.0949"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'0.0499"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").Value = "'18.53"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.38%  "
$ws.Range("D49").Value = "'0.0215"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("E50").Value = "  -5.38%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "
